$p = $ppt.ActivePresentation

# Slide 7 ("Objetivos Específicos") — fix typo in the last bullet of the
# content placeholder: "restes" -> "testes".
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$para = $tr.Paragraphs(5, 1)
$run = $para.Runs(1)

if ($run.Text -notlike "*restes*") {
    throw "Unexpected text in target paragraph: $($run.Text)"
}

$run.Text = $run.Text -replace "restes", "testes"
